$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("I2").Value = 0.7783932765807232
$ws.Range("J2").Value = 0.7783932765807231
$ws.Range("S2").Value = 0.7783932765807232
$ws.Range("T2").Value = 0.7783932765807231

# Row 3 updates
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.065288
$ws.Range("H3").Value = 0.195864
$ws.Range("I3").Value = 0.2216067234192769
$ws.Range("J3").Value = 0.2216067234192769
$ws.Range("Q3").Value = 0.0004146005626666666
$ws.Range("R3").Value = 0.003731405064
$ws.Range("S3").Value = 0.2216067234192769
$ws.Range("T3").Value = 0.2216067234192769
